$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 3
$ws.Range("H3").Value = 42499.5
$ws.Range("J3").Value = 42499.5
$ws.Range("L3").Value = 42499.5
$ws.Range("N3").Value = -42727.5
# Row 33
$ws.Range("H33").Value = 601.9167
$ws.Range("I33").Value = 302.33334
$ws.Range("K33").Value = 302.33334
$ws.Range("M33").Value = -73.33334000000002
# Row 53
$ws.Range("H53").Value = 1349.8334
$ws.Range("I53").Value = 1349.8334
$ws.Range("K53").Value = 1349.8334
$ws.Range("M53").Value = -712.8334
# Row 64
$ws.Range("H64").Value = 5201
$ws.Range("J64").Value = 5334.3335
$ws.Range("L64").Value = 5334.3335
$ws.Range("N64").Value = -5830.3335
# Row 67
$ws.Range("H67").Value = 5201
$ws.Range("J67").Value = 5334.3335
$ws.Range("L67").Value = 5334.3335
$ws.Range("N67").Value = -7050.3335
# Row 102
$ws.Range("H102").Value = 42499.5
$ws.Range("J102").Value = 42499.5
$ws.Range("L102").Value = 42499.5
$ws.Range("N102").Value = -48989.5
# Row 103
$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("M103").ClearContents()
# Row 105
$ws.Range("H105").Value = 50000
$ws.Range("J105").Value = 50000
$ws.Range("L105").Value = 50000
$ws.Range("N105").Value = -56988
# Row 111
$ws.Range("H111").Value = 1000
$ws.Range("I111").Value = 1000
$ws.Range("K111").Value = 3000
$ws.Range("M111").Value = 67
# Row 125
$ws.Range("H125").Value = 1299
$ws.Range("I125").Value = 1299
$ws.Range("K125").Value = 11691
$ws.Range("M125").Value = -9231
# Row 137
$ws.Range("H137").Value = 2307.2307
$ws.Range("I137").Value = 2974.625
$ws.Range("J137").Value = 1239.4
$ws.Range("K137").Value = 8923.875
$ws.Range("L137").Value = 3718.2
$ws.Range("M137").Value = -6373.875
$ws.Range("N137").Value = -8818.200000000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3544.4375
$ws.Range("I32").Value = 3515.0715
$ws.Range("K32").Value = 3515.0715
$ws.Range("M32").Value = -3228.0715
# Row 76
$ws.Range("H76").Value = 90000
$ws.Range("J76").Value = 90000
$ws.Range("L76").Value = 90000
$ws.Range("N76").Value = -90676
# Row 79
$ws.Range("H79").Value = 90000
$ws.Range("J79").Value = 90000
$ws.Range("L79").Value = 90000
$ws.Range("N79").Value = -92340
# Row 101
$ws.Range("H101").Value = 24166.5
$ws.Range("J101").Value = 24166.5
$ws.Range("L101").Value = 24166.5
$ws.Range("N101").Value = -30656.5
# Row 102
$ws.Range("H102").Value = 1497.5
$ws.Range("I102").Value = 2000
$ws.Range("K102").Value = 2000
$ws.Range("M102").Value = -378
# Row 119
$ws.Range("H119").Value = 49333.332
$ws.Range("J119").Value = 49333.332
$ws.Range("L119").Value = 49333.332
$ws.Range("N119").Value = -59009.332
# Row 122
$ws.Range("H122").Value = 21997.5
$ws.Range("I122").Value = 39998
$ws.Range("J122").Value = 3997
$ws.Range("K122").Value = 119994
$ws.Range("L122").Value = 11991
$ws.Range("M122").Value = -117544
$ws.Range("N122").Value = -16891

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()
# Row 99
$ws.Range("H99").Value = 1412.5
$ws.Range("I99").Value = 991.6667
$ws.Range("K99").Value = 991.6667
$ws.Range("M99").Value = 506.3333
# Row 134
$ws.Range("H134").Value = 763.6667
$ws.Range("I134").Value = 763.6667
$ws.Range("K134").Value = 2291.0001
$ws.Range("M134").Value = 243.9998999999998

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 7295.5
$ws.Range("I58").Value = 6238.75
$ws.Range("K58").Value = 6238.75
$ws.Range("M58").Value = -6035.75
# Row 62
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
# Row 65
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
# Row 132
$ws.Range("H132").Value = 6339.7334
$ws.Range("I132").Value = 3028
$ws.Range("K132").Value = 9084
$ws.Range("M132").Value = -6554
# Row 136
$ws.Range("H136").Value = 7295.5
$ws.Range("I136").Value = 6238.75
$ws.Range("K136").Value = 18716.25
$ws.Range("M136").Value = -16166.25

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 5300
$ws.Range("I4").Value = 5300
$ws.Range("K4").Value = 15900
$ws.Range("M4").Value = -15788

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 41
$ws.Range("H41").Value = 2567.8572
$ws.Range("I41").Value = 2495.8333
$ws.Range("J41").Value = 3000
$ws.Range("K41").Value = 2495.8333
$ws.Range("L41").Value = 3000
$ws.Range("M41").Value = -2140.8333
$ws.Range("N41").Value = -3710
# Row 101
$ws.Range("H101").Value = 44500
$ws.Range("J101").Value = 44500
$ws.Range("L101").Value = 44500
$ws.Range("N101").Value = -50990
# Row 113
$ws.Range("H113").Value = 4022.4666
$ws.Range("I113").Value = 4478.0835
$ws.Range("J113").Value = 2200
$ws.Range("K113").Value = 4478.0835
$ws.Range("L113").Value = 2200
$ws.Range("M113").Value = -2308.0835
$ws.Range("N113").Value = -6540
# Row 122
$ws.Range("H122").Value = 1793.5
$ws.Range("I122").Value = 1793.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5380.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2930.5
$ws.Range("N122").ClearContents()
# Row 127
$ws.Range("H127").Value = 99998.5
$ws.Range("J127").Value = 99998.5
$ws.Range("L127").Value = 99998.5
$ws.Range("N127").Value = -109918.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 374.375
$ws.Range("I22").Value = 374.375
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 374.375
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -79.375
$ws.Range("N22").ClearContents()
# Row 27
$ws.Range("H27").Value = 374.375
$ws.Range("I27").Value = 374.375
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 374.375
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -267.375
$ws.Range("N27").ClearContents()
# Row 55
$ws.Range("H55").Value = 374.4
$ws.Range("J55").Value = 285
$ws.Range("L55").Value = 285
$ws.Range("N55").Value = -631
# Row 68
$ws.Range("H68").Value = 2056.875
$ws.Range("I68").Value = 1992.5
$ws.Range("J68").Value = 2250
$ws.Range("K68").Value = 1992.5
$ws.Range("L68").Value = 2250
$ws.Range("M68").Value = -1243.5
$ws.Range("N68").Value = -3748
# Row 71
$ws.Range("H71").Value = 2056.875
$ws.Range("I71").Value = 1992.5
$ws.Range("J71").Value = 2250
$ws.Range("K71").Value = 9962.5
$ws.Range("L71").Value = 11250
$ws.Range("M71").Value = -6218.5
$ws.Range("N71").Value = -18738
# Row 82
$ws.Range("H82").Value = 1666.6666
$ws.Range("I82").Value = 1666.6666
$ws.Range("K82").Value = 1666.6666
$ws.Range("M82").Value = -1305.6666
# Row 85
$ws.Range("H85").Value = 1666.6666
$ws.Range("I85").Value = 1666.6666
$ws.Range("K85").Value = 1666.6666
$ws.Range("M85").Value = -418.6666
# Row 93
$ws.Range("H93").Value = 5249.6665
$ws.Range("I93").Value = 5249.6665
$ws.Range("K93").Value = 5249.6665
$ws.Range("M93").Value = -4001.6665
# Row 101
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
# Row 104
$ws.Range("H104").Value = 24411.334
$ws.Range("J104").Value = 24411.334
$ws.Range("L104").Value = 24411.334
$ws.Range("N104").Value = -31399.334
# Row 136
$ws.Range("H136").Value = 10400
$ws.Range("I136").Value = 10400
$ws.Range("K136").Value = 31200
$ws.Range("M136").Value = -28650

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 103
$ws.Range("H103").Value = 48000
$ws.Range("J103").Value = 48000
$ws.Range("L103").Value = 48000
$ws.Range("N103").Value = -50344
# Row 122
$ws.Range("H122").Value = 2429
$ws.Range("I122").Value = 1834.6666
$ws.Range("J122").Value = 2874.75
$ws.Range("K122").Value = 5503.9998
$ws.Range("L122").Value = 8624.25
$ws.Range("M122").Value = -3053.9998
$ws.Range("N122").Value = -13524.25
# Row 123
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
# Row 132
$ws.Range("H132").Value = 4000
$ws.Range("I132").Value = 4000
$ws.Range("K132").Value = 12000
$ws.Range("M132").Value = -9470
